$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "last updated" timestamp string (A1)
$ws.Range("A1").Value = "Datos actualizados a 27 de Marzo de 2020 a las 10:12"

# --- Step 1: update all country-name (column A) labels first, in one batch,
# so every shared-string text stays referenced throughout (avoids the shared
# string table dropping an entry that is still needed a few rows down).
$ws.Range("A36").Value = "Indonesia"
$ws.Range("A37").Value = "Rusia"
$ws.Range("A38").Value = "Rumania"
$ws.Range("A39").Value = "Arabia Saudita"
$ws.Range("A40").Value = "Finlandia"
$ws.Range("A41").Value = "Sudafrica"
$ws.Range("A43").Value = "Filipinas"
$ws.Range("A44").Value = "Islandia"
$ws.Range("A45").Value = "India"
$ws.Range("A46").Value = "Crucero"
$ws.Range("A56").Value = "Hong Kong"
$ws.Range("A57").Value = "Egipto"
$ws.Range("A58").Value = "Colombia"
$ws.Range("A59").Value = "Republica Dominicana"
$ws.Range("A60").Value = "Barein"
$ws.Range("A61").Value = "Serbia"

# --- Step 2: update the numeric statistics columns (B:H) for every changed row.
$ws.Range("B15").Value = 7129
$ws.Range("C15").Value = 220
$ws.Range("E15").Value = 6968
$ws.Range("B20").Value = 3380
$ws.Range("C20").Value = 8
$ws.Range("E20").Value = 3360
$ws.Range("B36").Value = 1046
$ws.Range("C36").Value = 153
$ws.Range("D36").Value = 46
$ws.Range("E36").Value = 913
$ws.Range("F36").Value = 0
$ws.Range("G36").Value = 9
$ws.Range("H36").Value = 87
$ws.Range("B37").Value = 1036
$ws.Range("C37").Value = 196
$ws.Range("D37").Value = 45
$ws.Range("E37").Value = 988
$ws.Range("F37").Value = 8
$ws.Range("G37").Value = 0
$ws.Range("H37").Value = 3
$ws.Range("B38").Value = 1029
$ws.Range("D38").Value = 94
$ws.Range("E38").Value = 911
$ws.Range("F38").Value = 29
$ws.Range("G38").Value = 1
$ws.Range("H38").Value = 24
$ws.Range("B39").Value = 1012
$ws.Range("D39").Value = 33
$ws.Range("E39").Value = 976
$ws.Range("F39").Value = 6
$ws.Range("H39").Value = 3
$ws.Range("B40").Value = 958
$ws.Range("D40").Value = 10
$ws.Range("E40").Value = 943
$ws.Range("F40").Value = 24
$ws.Range("G40").Value = 0
$ws.Range("H40").Value = 5
$ws.Range("B41").Value = 927
$ws.Range("D41").Value = 12
$ws.Range("E41").Value = 913
$ws.Range("F41").Value = 7
$ws.Range("G41").Value = 2
$ws.Range("H41").Value = 2
$ws.Range("B43").Value = 803
$ws.Range("C43").Value = 96
$ws.Range("D43").Value = 31
$ws.Range("F43").Value = 1
$ws.Range("G43").Value = 9
$ws.Range("H43").Value = 54
$ws.Range("B44").Value = 802
$ws.Range("C44").Value = 0
$ws.Range("D44").Value = 82
$ws.Range("E44").Value = 718
$ws.Range("F44").Value = 11
$ws.Range("H44").Value = 2
$ws.Range("B45").Value = 753
$ws.Range("C45").Value = 26
$ws.Range("D45").Value = 67
$ws.Range("E45").Value = 666
$ws.Range("F45").Value = 0
$ws.Range("H45").Value = 20
$ws.Range("B46").Value = 712
$ws.Range("D46").Value = 597
$ws.Range("E46").Value = 105
$ws.Range("F46").Value = 15
$ws.Range("H46").Value = 10
$ws.Range("B56").Value = 518
$ws.Range("C56").Value = 64
$ws.Range("D56").Value = 111
$ws.Range("E56").Value = 403
$ws.Range("F56").Value = 5
$ws.Range("H56").Value = 4
$ws.Range("B57").Value = 495
$ws.Range("D57").Value = 102
$ws.Range("E57").Value = 369
$ws.Range("H57").Value = 24
$ws.Range("B58").Value = 491
$ws.Range("D58").Value = 8
$ws.Range("E58").Value = 477
$ws.Range("H58").Value = 6
$ws.Range("B59").Value = 488
$ws.Range("D59").Value = 3
$ws.Range("E59").Value = 475
$ws.Range("F59").Value = 0
$ws.Range("H59").Value = 10
$ws.Range("B60").Value = 458
$ws.Range("D60").Value = 210
$ws.Range("E60").Value = 244
$ws.Range("F60").Value = 1
$ws.Range("H60").Value = 4
$ws.Range("B61").Value = 457
$ws.Range("D61").Value = 15
$ws.Range("E61").Value = 435
$ws.Range("F61").Value = 21
$ws.Range("H61").Value = 7
$ws.Range("E86").Value = 149
$ws.Range("G86").Value = 2
$ws.Range("H86").Value = 8
$ws.Range("B94").Value = 125
$ws.Range("C94").Value = 12
$ws.Range("E94").Value = 122
